$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Insert")

$values = 86,87,88,89,90,91,92,93,94,95
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

$ws.Activate()
$ws.Range("A2:A11").Select()
